# Scheduled runner update: refresh Market Board pricing + leve profit calcs
# across ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 225000
$ws.Range("J75").Value = 225000
$ws.Range("L75").Value = 225000
$ws.Range("N75").Value = -226872

$ws.Range("H76").Value = 3394.5
$ws.Range("I76").Value = 3090.5
$ws.Range("J76").Value = 3698.5
$ws.Range("K76").Value = 3090.5
$ws.Range("L76").Value = 3698.5
$ws.Range("M76").Value = -2775.5
$ws.Range("N76").Value = -4328.5

$ws.Range("H78").Value = 225000
$ws.Range("J78").Value = 225000
$ws.Range("L78").Value = 675000
$ws.Range("N78").Value = -684360

$ws.Range("H79").Value = 3394.5
$ws.Range("I79").Value = 3090.5
$ws.Range("J79").Value = 3698.5
$ws.Range("K79").Value = 3090.5
$ws.Range("L79").Value = 3698.5
$ws.Range("M79").Value = -1998.5
$ws.Range("N79").Value = -5882.5

$ws.Range("H86").Value = 74965.664
$ws.Range("I86").Value = 3698.5
$ws.Range("J86").Value = 217500
$ws.Range("K86").Value = 3698.5
$ws.Range("L86").Value = 217500
$ws.Range("M86").Value = -2575.5
$ws.Range("N86").Value = -219746

$ws.Range("H89").Value = 74965.664
$ws.Range("I89").Value = 3698.5
$ws.Range("J89").Value = 217500
$ws.Range("K89").Value = 18492.5
$ws.Range("L89").Value = 1087500
$ws.Range("M89").Value = -12876.5
$ws.Range("N89").Value = -1098732

$ws.Range("H106").Value = 11903.125
$ws.Range("I106").Value = 11460.714
$ws.Range("K106").Value = 11460.714
$ws.Range("M106").Value = -10829.714

$ws.Range("H112").Value = 2050.2104
$ws.Range("J112").Value = 2274.2307
$ws.Range("L112").Value = 6822.6921
$ws.Range("N112").Value = -9038.6921

$ws.Range("H132").Value = 2044013.6
$ws.Range("I132").Value = 2992.578
$ws.Range("K132").Value = 8977.734
$ws.Range("M132").Value = -6447.734

$ws.Range("H137").Value = 12686.105
$ws.Range("I137").Value = 13590.353
$ws.Range("J137").Value = 5000
$ws.Range("K137").Value = 40771.05899999999
$ws.Range("L137").Value = 15000
$ws.Range("M137").Value = -38221.05899999999
$ws.Range("N137").Value = -20100

$ws.Range("H141").Value = 6247.1665
$ws.Range("I141").Value = 6028.0625
$ws.Range("K141").Value = 18084.1875
$ws.Range("M141").Value = -12904.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7301.426
$ws.Range("I32").Value = 7322.635
$ws.Range("K32").Value = 7322.635
$ws.Range("M32").Value = -7035.635

$ws.Range("H48").Value = 174055
$ws.Range("J48").Value = 174055
$ws.Range("L48").Value = 174055
$ws.Range("N48").Value = -174823

$ws.Range("H69").Value = 213269.28
$ws.Range("J69").Value = 213269.28
$ws.Range("L69").Value = 213269.28
$ws.Range("N69").Value = -214767.28

$ws.Range("H72").Value = 213269.28
$ws.Range("J72").Value = 213269.28
$ws.Range("L72").Value = 639807.84
$ws.Range("N72").Value = -647295.84

$ws.Range("H102").Value = 9006.344999999999
$ws.Range("I102").Value = 11460.143
$ws.Range("K102").Value = 11460.143
$ws.Range("M102").Value = -9838.143

$ws.Range("H122").Value = 6636.625
$ws.Range("J122").Value = 8102.2
$ws.Range("L122").Value = 24306.6
$ws.Range("N122").Value = -29206.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 14623.27
$ws.Range("I99").Value = 17321.053
$ws.Range("K99").Value = 17321.053
$ws.Range("M99").Value = -15823.053

$ws.Range("H105").Value = 95962
$ws.Range("I105").Value = 168848.17
$ws.Range("K105").Value = 168848.17
$ws.Range("M105").Value = -167101.17

$ws.Range("H107").Value = 2125.1667
$ws.Range("I107").Value = 2243.5293
$ws.Range("K107").Value = 2243.5293
$ws.Range("M107").Value = -323.5293000000001

$ws.Range("H134").Value = 5188.8438
$ws.Range("I134").Value = 5740.5557
$ws.Range("J134").Value = 2209.6
$ws.Range("K134").Value = 17221.6671
$ws.Range("L134").Value = 6628.799999999999
$ws.Range("M134").Value = -14686.6671
$ws.Range("N134").Value = -11698.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H132").Value = 1855.0555
$ws.Range("I132").Value = 1905.6875
$ws.Range("K132").Value = 5717.0625
$ws.Range("M132").Value = -3187.0625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 231.05556
$ws.Range("I107").Value = 269.27274
$ws.Range("J107").Value = 171
$ws.Range("K107").Value = 269.27274
$ws.Range("L107").Value = 171
$ws.Range("M107").Value = 1650.72726
$ws.Range("N107").Value = -4011

$ws.Range("H126").Value = 13530
$ws.Range("I126").Value = 65006
$ws.Range("J126").Value = 7474
$ws.Range("K126").Value = 195018
$ws.Range("L126").Value = 22422
$ws.Range("M126").Value = -192548
$ws.Range("N126").Value = -27362

$ws.Range("H132").Value = 2250.612
$ws.Range("I132").Value = 2159.018
$ws.Range("K132").Value = 6477.054
$ws.Range("M132").Value = -3947.054

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 9204.5
$ws.Range("J22").Value = 1187.1428
$ws.Range("L22").Value = 1187.1428
$ws.Range("N22").Value = -1777.1428

$ws.Range("H27").Value = 9204.5
$ws.Range("J27").Value = 1187.1428
$ws.Range("L27").Value = 1187.1428
$ws.Range("N27").Value = -1401.1428

$ws.Range("H40").Value = 21415.033
$ws.Range("J40").Value = 14021.786
$ws.Range("L40").Value = 14021.786
$ws.Range("N40").Value = -14293.786

$ws.Range("H55").Value = 2180.182
$ws.Range("I55").Value = 426.16666
$ws.Range("J55").Value = 4285
$ws.Range("K55").Value = 426.16666
$ws.Range("L55").Value = 4285
$ws.Range("M55").Value = -253.16666
$ws.Range("N55").Value = -4631

$ws.Range("H61").Value = 8374.444
$ws.Range("J61").Value = 20000
$ws.Range("L61").Value = 20000
$ws.Range("N61").Value = -20404

$ws.Range("H68").Value = 3692.5881
$ws.Range("I68").Value = 2075.4167
$ws.Range("J68").Value = 7573.8
$ws.Range("K68").Value = 2075.4167
$ws.Range("L68").Value = 7573.8
$ws.Range("M68").Value = -1326.4167
$ws.Range("N68").Value = -9071.799999999999

$ws.Range("H71").Value = 3692.5881
$ws.Range("I71").Value = 2075.4167
$ws.Range("J71").Value = 7573.8
$ws.Range("K71").Value = 10377.0835
$ws.Range("L71").Value = 37869
$ws.Range("M71").Value = -6633.083500000001
$ws.Range("N71").Value = -45357

$ws.Range("H113").Value = 8374.444
$ws.Range("J113").Value = 20000
$ws.Range("L113").Value = 20000
$ws.Range("N113").Value = -24340

$ws.Range("H122").Value = 5523.316
$ws.Range("I122").Value = 5584.8823
$ws.Range("K122").Value = 16754.6469
$ws.Range("M122").Value = -14304.6469

$ws.Range("H136").Value = 4356.143
$ws.Range("I136").Value = 3320.3
$ws.Range("K136").Value = 9960.900000000001
$ws.Range("M136").Value = -7410.900000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5256.4653
$ws.Range("I122").Value = 2459.8635
$ws.Range("J122").Value = 8186.2383
$ws.Range("K122").Value = 7379.5905
$ws.Range("L122").Value = 24558.7149
$ws.Range("M122").Value = -4929.5905
$ws.Range("N122").Value = -29458.7149

$ws.Range("H125").Value = 150000
$ws.Range("J125").Value = 150000
$ws.Range("L125").Value = 150000
$ws.Range("N125").Value = -159840

$ws.Range("H132").Value = 9823.352999999999
$ws.Range("I132").Value = 12875.686
$ws.Range("K132").Value = 38627.058
$ws.Range("M132").Value = -36097.058

$ws.Range("H136").Value = 345212.62
$ws.Range("I136").Value = 418617.06
$ws.Range("K136").Value = 1255851.18
$ws.Range("M136").Value = -1253301.18
